# Add a "2022-Q4" fund-holding sheet (copied from "2022-Q3" then re-populated
# with the new quarter's figures) and record it in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet right before "2022-Q3" by
#    duplicating "2022-Q3" (keeps headers/column widths/styles intact)
#    and then overwriting its data rows with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3, $null)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

$q4rows = @(
  @{r=2; B="009225"; C="天弘中证中美互联网指数（QDII）A";         D="1.27"; E="94.90"; F="4.65"; G="0.0591"; H=8},
  @{r=3; B="012751"; C="建信纳斯达克100指数（QDII）A 美元现汇"; D="1.06"; E="82.28"; F="3.50"; G="0.0371"; H=5},
  @{r=4; B="012752"; C="建信纳斯达克100指数（QDII）C 人民币";     D="1.06"; E="82.28"; F="3.50"; G="0.0371"; H=5},
  @{r=5; B="012753"; C="建信纳斯达克100指数（QDII）C 美元现汇"; D="1.06"; E="82.28"; F="3.50"; G="0.0371"; H=5},
  @{r=6; B="009226"; C="天弘中证中美互联网指数（QDII）C";         D="0.63"; E="94.90"; F="4.65"; G="0.0293"; H=8}
)

foreach ($row in $q4rows) {
    $r = $row.r
    # Fund code / decimal-looking figures must stay TEXT (leading zeros,
    # trailing zeros matter) -- a leading apostrophe forces text entry,
    # then ClearFormats() drops the resulting "number stored as text"
    # quote-prefix flag so the cell style matches the rest of the sheet.
    $wsQ4.Range("B$r").Value = "'" + $row.B
    $wsQ4.Range("C$r").Value = $row.C
    $wsQ4.Range("D$r").Value = "'" + $row.D
    $wsQ4.Range("E$r").Value = "'" + $row.E
    $wsQ4.Range("F$r").Value = "'" + $row.F
    $wsQ4.Range("G$r").Value = "'" + $row.G
    $wsQ4.Range("H$r").Value = $row.H

    $wsQ4.Range("B$r").ClearFormats()
    $wsQ4.Range("D$r`:G$r").ClearFormats()
}

# ---------------------------------------------------------------------
# 2. Insert a new row into "总计" for the 2022-Q4 summary, pushing the
#    existing 2022-Q3 / 2021-Q4 / 2021-Q2 rows down by one.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Copy row-index column formatting down onto the newly inserted A2.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.2

# Keep the running row-index column sequential after the insert.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
